$d = $word.ActiveDocument

# Locate the paragraph that ends with the "main constraint ... 100 and 1000."
# sentence (item 2 "Break the problem apart") - this is where the new
# "3) Identify potential solutions" content needs to be inserted right after,
# and before the trailing blank paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "every 10 units the count stops in a different finger") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the 'Break the problem apart' constraint paragraph to insert after."
}

$p39 = $d.Paragraphs.Item($targetIndex)

# The trailing "_GoBack" bookmark currently sits at the end of this
# paragraph; it needs to end up after the newly-inserted content instead,
# so remove it now and re-add it in the right spot below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insertion point: right before the paragraph mark of the target paragraph
# (i.e. at the very end of its visible text).
$insertionPoint = $d.Range($p39.Range.End - 1, $p39.Range.End - 1)

$newContentXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="left" w:pos="220"/><w:tab w:val="left" w:pos="720"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>3)  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Identify potential solutions </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="left" w:pos="220"/><w:tab w:val="left" w:pos="720"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:cs="Calibri"/><w:bCs/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Calibri"/><w:bCs/><w:szCs w:val="32"/></w:rPr><w:t>I count from 1 to 100 stopping in 10, 20, 30, 40, 5</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Calibri"/><w:bCs/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">0, 60, 70, 80, 90, 100, and realized that the count stops in the first finger twice every two times and the same with the ring finger. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$null = $insertionPoint.InsertXML($newContentXml)
